$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 3 de Agosto de 2020 a las 00:54"

# Colombia overtook Iran in ranking: row 13 becomes Colombia (new data),
# row 14 becomes Iran (using what used to be Iran's figures at row 13)
$ws.Range("A13").Value = "Colombia"
$ws.Range("A14").Value = "Iran"

# Updated case numbers for several countries
$ws.Range("B4").Value = 4811959
$ws.Range("C4").Value = 47350
$ws.Range("D4").Value = 2376567
$ws.Range("E4").Value = 2277070
$ws.Range("G4").Value = 424
$ws.Range("H4").Value = 158322

$ws.Range("B10").Value = 428850
$ws.Range("C10").Value = 6667
$ws.Range("D10").Value = 294187
$ws.Range("E10").Value = 115049
$ws.Range("G10").Value = 206
$ws.Range("H10").Value = 19614

$ws.Range("B13").Value = 317651
$ws.Range("C13").Value = 11470
$ws.Range("D13").Value = 167239
$ws.Range("E13").Value = 139762
$ws.Range("G13").Value = 320
$ws.Range("H13").Value = 10650

$ws.Range("B14").Value = 309437
$ws.Range("C14").Value = 2685
$ws.Range("D14").Value = 268102
$ws.Range("E14").Value = 24145
$ws.Range("G14").Value = 208
$ws.Range("H14").Value = 17190

$ws.Range("B42").Value = 67453
$ws.Range("C42").Value = 1070
$ws.Range("D42").Value = 41038
$ws.Range("E42").Value = 24944
$ws.Range("G42").Value = 22
$ws.Range("H42").Value = 1471

$ws.Range("B50").Value = 43841
$ws.Range("C50").Value = 304
$ws.Range("D50").Value = 20308
$ws.Range("E50").Value = 22645
$ws.Range("G50").Value = 5
$ws.Range("H50").Value = 888

$ws.Range("B57").Value = 36689
$ws.Range("C57").Value = 853
$ws.Range("D57").Value = 25906
$ws.Range("E57").Value = 9772

$ws.Range("D58").Value = 31500
$ws.Range("E58").Value = 2069

$ws.Range("B69").Value = 20206
$ws.Range("C69").Value = 763
$ws.Range("D69").Value = 11404
$ws.Range("E69").Value = 8628
$ws.Range("G69").Value = 5
$ws.Range("H69").Value = 174

$ws.Range("B75").Value = 16800
$ws.Range("C75").Value = 101
$ws.Range("D75").Value = 11605
$ws.Range("E75").Value = 4811

$ws.Range("B81").Value = 11955
$ws.Range("C81").Value = 119
$ws.Range("D81").Value = 6420
$ws.Range("E81").Value = 5147
$ws.Range("G81").Value = 3
$ws.Range("H81").Value = 388

$ws.Range("B86").Value = 9268
$ws.Range("C86").Value = 15
$ws.Range("E86").Value = 261

$ws.Range("B114").Value = 3220
$ws.Range("C114").Value = 8
$ws.Range("D114").Value = 1598
$ws.Range("E114").Value = 1529

$ws.Range("B122").Value = 2541
$ws.Range("C122").Value = 6
$ws.Range("D122").Value = 1943
$ws.Range("E122").Value = 474

$ws.Range("B139").Value = 1286
$ws.Range("C139").Value = 8
$ws.Range("D139").Value = 1011
$ws.Range("G139").Value = 1
$ws.Range("H139").Value = 36

